$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("I2").Value = 0.415
$ws.Range("J2").Value = 23.166
$ws.Range("L2").Value = 98.25177359
$ws.Range("M2").Value = 0.037469988

$ws.Range("I3").Value = 1.715
$ws.Range("J3").Value = 24.82
$ws.Range("L3").Value = 93.13621678
$ws.Range("M3").Value = 0.294448637

$ws.Range("I4").Value = 1.495
$ws.Range("J4").Value = 11.792
$ws.Range("L4").Value = 85.87427554
$ws.Range("M4").Value = 6.172127348

$ws.Range("I5").Value = 0.944
$ws.Range("J5").Value = 23.333
$ws.Range("L5").Value = 95.99713334
$ws.Range("M5").Value = 0.071979164

$ws.Range("I6").Value = 1.405
$ws.Range("J6").Value = 20.706
$ws.Range("L6").Value = 93.37792351
$ws.Range("M6").Value = 2.714900366

$ws.Range("I6:M6").Select()
